# Update the "monthly_gross_earnings_effect_*" labels in column A (rows 2-13)
# to "yearly_gross_earnings_effect_*", and move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 -replace "^monthly_", "yearly_"
}

$ws.Range("A27").Select()
